# Update the "fs_results" feature-sorting table (Sheet_name_1) with the
# final values used for testing of the dp4 dissertation.
# Column A (index) values are left untouched; several cells in rows
# 5, 8, 10 and 11 are cleared entirely (no replacement value), matching
# the source data that no longer has those measurements.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -35
$ws.Range("C2").Value = 7
$ws.Range("D2").Value = 6
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 19
$ws.Range("G2").Value = 9
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = -2

$ws.Range("B3").Value = -24
$ws.Range("C3").Value = 13
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = -48
$ws.Range("G3").Value = 120
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = -4

$ws.Range("B4").Value = 89
$ws.Range("C4").Value = -96
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = -29
$ws.Range("H4").Value = -1
$ws.Range("I4").Value = 94

$ws.Range("B5").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = -32
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 1

$ws.Range("B6").Value = -16
$ws.Range("C6").Value = 8
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = 13
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 107

$ws.Range("B7").Value = -10
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = 76
$ws.Range("E7").Value = 13
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = -37
$ws.Range("H7").Value = -3
$ws.Range("I7").Value = -3

$ws.Range("B8").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("D8").Value = 79
$ws.Range("E8").Value = 15
$ws.Range("F8").Value = 3
$ws.Range("G8").Value = -40
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()

$ws.Range("B9").Value = -12
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 85
$ws.Range("E9").Value = 56
$ws.Range("F9").Value = -13
$ws.Range("G9").Value = 9
$ws.Range("H9").Value = 6
$ws.Range("I9").Value = -4

$ws.Range("B10").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("D10").Value = 112
$ws.Range("E10").Value = -8
$ws.Range("F10").Value = -11
$ws.Range("G10").Value = -4
$ws.Range("H10").ClearContents()
$ws.Range("I10").ClearContents()

$ws.Range("B11").Value = -10
$ws.Range("C11").Value = 4
$ws.Range("D11").Value = 115
$ws.Range("E11").Value = 18
$ws.Range("F11").Value = 5
$ws.Range("G11").Value = -22
$ws.Range("H11").ClearContents()
$ws.Range("I11").ClearContents()

